$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TMA map")

# --- Row 2 updates (right block shifts up) ---
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 6
$ws.Range("K2").Value = 7
$ws.Range("L2").Value = 7

# --- Row 3 updates: left block collapses to a single pair, right block shifts ---
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 3
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9
# K3 / L3 stay 8

# --- Row 4 is new: insert data (both left and right blocks) ---
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 5
$ws.Range("I4").Value = 10
$ws.Range("J4").Value = 10
$ws.Range("K4").Value = 11
$ws.Range("L4").Value = 11

# --- Row 7 updates ---
$ws.Range("C7").Value = 12
$ws.Range("D7").Value = 12
$ws.Range("E7").Value = 13
$ws.Range("F7").Value = 13
$ws.Range("I7").Value = 16
$ws.Range("J7").Value = 16
$ws.Range("K7").Value = 17
$ws.Range("L7").Value = 17

# --- Row 8 updates ---
$ws.Range("C8").Value = 15
$ws.Range("D8").Value = 15
$ws.Range("E8").Value = 14
$ws.Range("F8").Value = 14
$ws.Range("I8").Value = 19
$ws.Range("J8").Value = 19
$ws.Range("K8").Value = 18
$ws.Range("L8").Value = 18

# --- Update selection to reflect final cursor position ---
$ws.Range("D10").Select()
